$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-12 hold the existing match-log data. The update duplicates that
# data (in a specific re-shuffled row order) and appends it as rows 13-23,
# growing the sheet from A1:K12 to A1:K23.
$srcRows = @(11, 4, 10, 6, 8, 12, 9, 2, 5, 7, 3)
$firstDestRow = 13
$lastDestRow = $firstDestRow + $srcRows.Count - 1

# Columns G:K (totalRuns, totalBalls, total4s, total6s, sr) hold numbers
# that are stored as text in the source sheet (e.g. "7", "53.84"). Format
# that block as Text up front so the values we copy in stay text instead
# of being auto-coerced to numbers.
$ws.Range("G$firstDestRow`:K$lastDestRow").NumberFormat = "@"

$destRow = $firstDestRow
foreach ($srcRow in $srcRows) {
    for ($col = 1; $col -le 11; $col++) {
        $srcCell = $ws.Cells.Item($srcRow, $col)
        $destCell = $ws.Cells.Item($destRow, $col)
        $destCell.Value = $srcCell.Value2
    }
    $destRow++
}
